$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same style as A1/A2 to A3 (resolves type/formatting for the row)
$ws.Range("A3").Value = $ws.Range("A3").Value
$ws.Range("A1").Copy()
$ws.Range("A3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Ensure the text value for A3 remains "jaipur J"
$ws.Range("A3").Value = "jaipur J"

# Add new row with the new city
$ws.Range("A4").Value = "Gurgain"

# Update selection to the next empty cell, matching Excel's default behavior after data entry
$ws.Range("A5").Select() | Out-Null
